$d = $word.ActiveDocument

# 1. Contract date (paragraph 3: "Puerto Peñasco, Sonora. 27/septiembre/2017")
$d.Paragraphs(3).Range.Find.Execute("27/septiembre/2017", $true, $false, $false, $false, $false,
                         $true, 1, $false, "30/octubre/2017", 2)

# 2. Client name, first occurrence (paragraph 5, intro clause)
$d.Paragraphs(5).Range.Find.Execute("Ofelia  Martinez Zamora", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Cliente De Prueba Tester", 2)

# 3. Event date (paragraph 6)
$d.Paragraphs(6).Range.Find.Execute("28/septiembre/2017", $true, $false, $false, $false, $false,
                         $true, 1, $false, "08/noviembre/2017", 2)

# 4. Event time (paragraph 6)
$d.Paragraphs(6).Range.Find.Execute("00:00hrs.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "14:00hrs.", 2)

# 5. Number of people (paragraph 7: "Número de Personas: 20")
$d.Paragraphs(7).Range.Find.Execute("20", $true, $false, $false, $false, $false,
                         $true, 1, $false, "8", 2)

# 6. Place description (paragraph 8)
$d.Paragraphs(8).Range.Find.Execute("una descripcion", $true, $false, $false, $false, $false,
                         $true, 1, $false, "LAs Palomas Resort", 2)

# 7. Dish (paragraph 9)
$d.Paragraphs(9).Range.Find.Execute("PAto con pollo, 2 tiempos", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Chao Ming Pollo, 3 tiempos", 2)

# 8. Services included (paragraph 11)
$d.Paragraphs(11).Range.Find.Execute("2 Platos, 1 Vasos High Ball", $true, $false, $false, $false, $false,
                         $true, 1, $false, "4 Copas, 1 Barra de Licores y Cerveza, 2 Vasos High Ball, 1 Mesa de Postres, 1 Bocadillos, 1 Meseros", 2)

# 9. Price (paragraph 13)
$d.Paragraphs(13).Range.Find.Execute("$13,500.00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "$20,041.00", 2)

# 10. Advance payment (paragraph 13)
$d.Paragraphs(13).Range.Find.Execute("$0.00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "$1,500.00", 2)

# 11. Remove "Solo por tener un valor" paragraph entirely (paragraph 15)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`n`a") -eq "Solo por tener un valor") {
        $p.Range.Delete()
        break
    }
}

# 12. Client name, second occurrence (signature block)
$d.Content.Find.Execute("Ofelia  Martinez Zamora", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Cliente De Prueba Tester", 2)
